# Update title and licensing slides for ATPESC
#
# This script applies the content edits from the commit:
#   - Handout master / notes master "Update automatically" date
#     placeholders move from 8/7/2021 to 8/10/2021.
#   - Slide 2 ("License, Citation and Acknowledgements") citation
#     paragraph is re-worded for the ATPESC track, and the DOI is
#     updated to the new figshare record.

$p = $ppt.ActivePresentation

# --- 1) Handout master date placeholder -----------------------------
$handoutDate = $p.HandoutMaster.HeadersFooters.DateAndTime
$handoutDate.Text = "8/10/2021"

# --- 2) Notes master date placeholder --------------------------------
$notesDate = $p.NotesMaster.HeadersFooters.DateAndTime
$notesDate.Text = "8/10/2021"

# --- 3) Slide 2: citation + DOI text ---------------------------------
$slide = $p.Slides.Item(2)
$contentShape = $slide.Shapes.Item(2)
$textRange = $contentShape.TextFrame.TextRange

$oldCitation = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Patricia A. Grubel, Rinku K. Gupta, and David M. Rogers, Better Scientific Software tutorial, in ISC High Performance, online, 2021. DOI: "
$newCitation = "The requested citation the overall tutorial is: David E. Bernholdt, Anshu Dubey, Rinku K. Gupta, and David M. Rogers, Software Productivity and Sustainability track, in Argonne Training Program on Extreme-Scale Computing (ATPESC), online, 2021. DOI: "

$fullText = $textRange.Text
$startIdx = $fullText.IndexOf($oldCitation)
if ($startIdx -ge 0) {
    $run = $textRange.Characters($startIdx + 1, $oldCitation.Length)
    $run.Text = $newCitation
}

$oldDoi = "10.6084/m9.figshare.14642520"
$newDoi = "10.6084/m9.figshare.15130590"

$fullText2 = $textRange.Text
$doiIdx = $fullText2.IndexOf($oldDoi)
if ($doiIdx -ge 0) {
    $doiRun = $textRange.Characters($doiIdx + 1, $oldDoi.Length)
    $doiRun.Text = $newDoi
}
